# Method of Hierarchy Analysis - formatting pass
# - title cell A1
# - re-balance column A width
# - distinguish label cells (borders + vertical/wrap only) from the
#   pairwise-comparison matrix cells (which also get horizontal centering)
# - the "spacer" cells under the merged note block lose their old
#   vertical/wrap formatting, keeping only the font + border

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title text for the top-left corner of the matrix
$ws.Range("A1").Value = "Стоимость ТС"

# Column A a bit wider to fit the title
$ws.Columns.Item(1).ColumnWidth = 16

# --- Label-style cells: header row, row labels, summary notes (H8:H10) ---
# font + thin box border + vertical-center + wrap, no horizontal centering
$rngLabel = $ws.Range("A1:H1,A2:A7,H8:H10")
$rngLabel.Font.Name = "Times New Roman"
$rngLabel.Font.Size = 12
$rngLabel.Borders.LineStyle = 1
$rngLabel.HorizontalAlignment = 1
$rngLabel.VerticalAlignment = -4108
$rngLabel.WrapText = $true

# --- Data-matrix cells: pairwise comparison values, centered both ways ---
$rngData = $ws.Range("B2:H7")
$rngData.Font.Name = "Times New Roman"
$rngData.Font.Size = 12
$rngData.Borders.LineStyle = 1
$rngData.HorizontalAlignment = -4108
$rngData.VerticalAlignment = -4108
$rngData.WrapText = $true

# --- Spacer cells below the matrix (under the merged note block) ---
# keep the font + border, but drop the inherited vertical-center/wrap
$rngSpacer = $ws.Range("A8:G10")
$rngSpacer.Font.Name = "Times New Roman"
$rngSpacer.Font.Size = 12
$rngSpacer.Borders.LineStyle = 1
$rngSpacer.HorizontalAlignment = 1
$rngSpacer.VerticalAlignment = -4107
$rngSpacer.WrapText = $false
